$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 updates
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: clear B2, D2, E2; update C2
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = 5.3319794989134781
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3 updates
$ws.Range("B3").Value = 5.6375100864256718
$ws.Range("C3").Value = 6.0599002366774419
$ws.Range("D3").Value = 7.1032736555109457
$ws.Range("E3").Value = 3.0162978350576486

# Update selection to match new used range B1:E3
$ws.Range("B1:E3").Select()
